$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sreyas")

# New rows of data appended to the bottom of the table (rows 21-23)
$data = @(
    @("india", "sri lanka", 777,  90, 45344, 500),
    @("barca", "real",      6655, 90, 45344, 500),
    @("ooo",   "pppp",      1245, 90, 45344, 500)
)

$startRow = 21
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).NumberFormat = $ws.Cells.Item($r - 1, 5).NumberFormat
    $ws.Cells.Item($r, 6).Value = $row[5]
}
